# Update gh-pages to output generated at 456a3b4
# Increment "想去人数" (want-to-go count) values on both the "展览" and
# "全部类型" sheets: F2 523 -> 524, F4 6 -> 7.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 524
    $ws.Range("F4").Value = 7
}
